# Adds three more Knowledge-basis test rows (branch_and_bound predicate tests)
# to the Prolog Sprint1 test sheet, plus a note row, matching the commit:
# "More tests of our predicates added to the file of the tests"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: branch_and_bound('Vicky','JoseCid',L) ---
$ws.Range("B20").Value = "Knowledge basis"
$ws.Range("C20").Value = "branch_and_bound('Vicky','JoseCid',L)"
$ws.Range("D20").Value = "L=['Vicky','Joao','Tiago','Diogo','Francisco','JoseCid']"
$ws.Range("E20").Value = "False"

# --- Row 21: branch_and_bound('Simao','Maria',L) ---
$ws.Range("B21").Value = "Knowledge basis"
$ws.Range("C21").Value = "branch_and_bound('Simao','Maria',L)"
$ws.Range("D21").Value = "L=[Simao','Artur','Tiago','Stephanie','Maria']"
$ws.Range("E21").Value = "L=[Simao','Artur','Tiago','Stephanie','Maria']"

# --- Row 30: add a note in column B (next to the already-present C30 cell) ---
# (shared string allocated here so it lands before the next new string, to
# match the original author's edit order)
$ws.Range("B30").Value = "Note: Some bugs of branch_and_bound need to be fixed"

# --- Row 22: branch_and_bound('Joao','Tiago',L) ---
$ws.Range("B22").Value = "Knowledge basis"
$ws.Range("C22").Value = "branch_and_bound('Joao','Tiago',L)"
$ws.Range("D22").Value = "L=['Joao','Simao','Artur',Tiago']"
$ws.Range("E22").Value = "L=['Joao','Simao','Artur',Tiago']"
$ws.Range("E22").Font.Underline = $true

# --- Update the view's selection to match the edited area ---
$ws.Range("E22").Select()
